# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G on the active sheet holds the pitcher's per-game "K" (strikeout) count.
# The data was regenerated from the source stats (K instead of Strike#), so the
# values in G2:G70 are rewritten with their recalculated counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..70 (column G), in row order.
$kValues = @(
    0, 1, 0, 0, 0, 1, 1, 1, 1, 1,
    2, 2, 0, 1, 2, 3, 1, 2, 0, 1,
    1, 0, 2, 0, 0, 4, 0, 3, 1, 1,
    0, 1, 3, 0, 1, 3, 2, 0, 2, 0,
    1, 3, 1, 0, 1, 1, 1, 0, 1, 1,
    4, 1, 4, 1, 2, 0, 0, 2, 1, 4,
    2, 5, 1, 4, 3, 4, 7, 3, 3
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
